$d = $word.ActiveDocument

$searchStart = 0
while ($true) {
    $rng = $d.Range($searchStart, $d.Content.End)
    $found = $rng.Find.Execute("<id>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        break
    }

    $tagStart = $rng.Start
    $tagEnd = $rng.End

    # Peek at the text following "<id>" to find the matching "</id>" and
    # recover the id value in between. Grow the lookahead window in case
    # the id value turns out to be unexpectedly long.
    $docEnd = $d.Content.End
    $window = 200
    $idxClose = -1
    $peek = ""
    while ($idxClose -lt 0) {
        $peekEnd = [Math]::Min($tagEnd + $window, $docEnd)
        $peek = $d.Range($tagEnd, $peekEnd).Text
        $idxClose = $peek.IndexOf("</id>")
        if ($idxClose -ge 0) {
            break
        }
        if ($peekEnd -ge $docEnd) {
            break
        }
        $window = $window * 4
    }
    if ($idxClose -lt 0) {
        # Not the pattern we expect; skip past this tag and continue.
        $searchStart = $tagEnd
        continue
    }

    $idValue = $peek.Substring(0, $idxClose)
    $fullText = "<id>" + $idValue + "</id>"

    # Overwrite the "<id>" run's text with the full merged string - it
    # takes on that run's own formatting (Courier New / 7f6000 / 18).
    $mergeRange = $d.Range($tagStart, $tagEnd)
    $mergeRange.Text = $fullText

    # The old "idValue" + "</id>" runs now immediately follow; remove them.
    $staleStart = $tagStart + $fullText.Length
    $staleEnd = $staleStart + $idValue.Length + 5  # 5 == Len("</id>")
    $staleRange = $d.Range($staleStart, $staleEnd)
    $staleRange.Text = ""

    $searchStart = $staleStart
}
